$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

$ws.Range("B4").Value = 0.66666666666666663
$ws.Range("B5").Value = 0.91891891891891897
$ws.Range("B6").Value = 0.69387755102040816
$ws.Range("B7").Value = 0.79069767441860461
